$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("survey")
$ws.Range("A1:M16").ClearContents()
$ws.Range("A1").Value = "clause"
$ws.Range("B1").Value = "condition"
$ws.Range("C1").Value = "type"
$ws.Range("D1").Value = "values_list"
$ws.Range("E1").Value = "name"
$ws.Range("F1").Value = "display.prompt"
$ws.Range("G1").Value = "display.hint.text"
$ws.Range("H1").Value = "inputAttributes.min"
$ws.Range("I1").Value = "inputAttributes.step"
$ws.Range("J1").Value = "calculation"
$ws.Range("K1").Value = "constraint"
$ws.Range("L1").Value = "display.constraint_message"
$ws.Range("M1").Value = "hideInContents"
$ws.Range("A2").Value = "begin screen"
$ws.Range("C3").Value = "text"
$ws.Range("E3").Value = "hh_death_id"
$ws.Range("F3").Value = "q65a"
$ws.Range("G3").Value = "ABC-123-700"
$ws.Range("K3").Value = "/^[A-Z]{3}-[0-9]{3}-7[0-9]{2}`$/.test(data('hh_death_id'))"
$ws.Range("L3").Value = "extid_format"
$ws.Range("C4").Value = "text"
$ws.Range("E4").Value = "hh_death_name"
$ws.Range("F4").Value = "q65b"
$ws.Range("K4").Value = "!/\p{N}/u.test(data('hh_death_name'))"
$ws.Range("L4").Value = "name_number"
$ws.Range("C5").Value = "text"
$ws.Range("E5").Value = "hh_death_surname"
$ws.Range("F5").Value = "q65c"
$ws.Range("K5").Value = "!/\p{N}/u.test(data('hh_death_surname'))"
$ws.Range("L5").Value = "surname_number"
$ws.Range("C6").Value = "select_one"
$ws.Range("D6").Value = "gender"
$ws.Range("E6").Value = "hh_death_gender"
$ws.Range("F6").Value = "q65d"
$ws.Range("C7").Value = "note"
$ws.Range("F7").Value = "q65e"
$ws.Range("A8").Value = "if"
$ws.Range("B8").Value = "not(selected(data('hh_death_date_dk'), 'dk'))"
$ws.Range("C9").Value = "birth_date"
$ws.Range("E9").Value = "hh_death_date"
$ws.Range("K9").Value = "selected(data('hh_death_date_dk'), 'dk') || !data('hh_death_date') || (!data('hh_death_date').isBefore('2020-03-01') && !data('hh_death_date').isAfter())"
$ws.Range("L9").Value = "invalid_death_date"
$ws.Range("M9").Value = 1
$ws.Range("A10").Value = "end if"
$ws.Range("C11").Value = "select_multiple"
$ws.Range("D11").Value = "dk"
$ws.Range("E11").Value = "hh_death_date_dk"
$ws.Range("M11").Value = 1
$ws.Range("A12").Value = "if"
$ws.Range("B12").Value = "selected(data('hh_death_date_dk'), 'dk')"
$ws.Range("C13").Value = "assign"
$ws.Range("E13").Value = "hh_death_date"
$ws.Range("J13").Value = "null"
$ws.Range("A14").Value = "end if"
$ws.Range("C15").Value = "integer"
$ws.Range("E15").Value = "hh_death_age"
$ws.Range("F15").Value = "q65f"
$ws.Range("H15").Value = 0
$ws.Range("I15").Value = 1
$ws.Range("A16").Value = "end screen"
$ws.Range("F11").Select()

$ws6 = $wb.Worksheets.Item("table_specific_translations")
$ws6.Range("A14").Value = "invalid_death_date"
$ws6.Range("B14").Value = "Date cannot be in the future<br>Date should be in the past 12 months"
$ws6.Range("B14").Select()

$ws4 = $wb.Worksheets.Item("model")
$ws4.Range("B8").Select()

$ws5 = $wb.Worksheets.Item("settings")
$ws5.Range("B6").Select()

$ws.Activate()
